$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) updates - these are stored as text in the sheet, so a
# leading apostrophe forces Excel to keep them as text instead of coercing
# them to numbers.
$ws.Range("D2").Formula  = "'244.79"
$ws.Range("D3").Formula  = "'23.05"
$ws.Range("D4").Formula  = "'5.410"
$ws.Range("D5").Formula  = "'0.06018"
$ws.Range("D6").Formula  = "'3.395"
$ws.Range("D7").Formula  = "'0.8094"
$ws.Range("D8").Formula  = "'0.9267"
$ws.Range("D9").Formula  = "'0.1424"
$ws.Range("D10").Formula = "'0.07439"
$ws.Range("D11").Formula = "'0.03376"
$ws.Range("D12").Formula = "'0.03042"
$ws.Range("D13").Formula = "'0.09361"
$ws.Range("D14").Formula = "'3.933"
$ws.Range("D15").Formula = "'0.001587"
$ws.Range("D16").Formula = "'0.04825"
$ws.Range("D17").Formula = "'0.0005942"
$ws.Range("D18").Formula = "'0.005433"
$ws.Range("D19").Formula = "'0.004148"
$ws.Range("D20").Formula = "'0.0009868"
$ws.Range("D21").Formula = "'0.00008703"
$ws.Range("D23").Formula = "'6.436"
$ws.Range("D24").Formula = "'2.185"
$ws.Range("D27").Formula = "'0.0002447"
$ws.Range("D40").Formula = "'0.03966"
$ws.Range("D41").Formula = "'0.006412"
$ws.Range("D43").Formula = "'0.002901"
$ws.Range("D44").Formula = "'0.006559"
$ws.Range("D45").Formula = "'0.00005208"
$ws.Range("D46").Formula = "'0.00000000750"
$ws.Range("D47").Formula = "'0.0005802"
$ws.Range("D48").Formula = "'0.8503"
$ws.Range("D49").Formula = "'0.002274"
$ws.Range("D50").Formula = "'0.00002101"

# Volume(1h) label (column E) updates - plain text, "Worstin24h" suffix
# moved from the One (row 17) coin to the Bolo (row 49) coin.
$ws.Range("E17").Value = "16OneONE"
$ws.Range("E49").Value = "48BOLOBOLOWorstin24h"
